# Group 4 Assignment 4 (final version) - reorder rate rows within each
# "year" block (the table was re-sorted, flipping the order of the
# rate_value/year rows inside each 2006 / 2007 group).
#
# 2006 block: rows 2-4  -> reverse order (row 2 <-> row 4, row 3 unchanged)
# 2007 block: rows 5-9  -> reverse order (row 5 <-> row 9, row 6 <-> row 8,
#                           row 7 unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2) {
    $a1 = $ws.Cells.Item($r1, 1).Value2
    $e1 = $ws.Cells.Item($r1, 5).Value2
    $a2 = $ws.Cells.Item($r2, 1).Value2
    $e2 = $ws.Cells.Item($r2, 5).Value2

    $ws.Cells.Item($r1, 1).Value = $a2
    $ws.Cells.Item($r1, 5).Value = $e2
    $ws.Cells.Item($r2, 1).Value = $a1
    $ws.Cells.Item($r2, 5).Value = $e1
}

# 2006 block (rows 2-4)
Swap-Rows $ws 2 4

# 2007 block (rows 5-9)
Swap-Rows $ws 5 9
Swap-Rows $ws 6 8
